$d = $word.ActiveDocument

$targets = @(
    "Folders with more in them are rendered as visibly larger wormholes",
    "Will have to implement a secondary traversal algorithm"
)

foreach ($para in $d.Paragraphs) {
    $text = $para.Range.Text
    foreach ($t in $targets) {
        if ($text -like "*$t*") {
            $para.Range.Font.StrikeThrough = $true
        }
    }
}
